# Insert a new column "gross_e5_tu" right after column G (gross_e5),
# pushing the existing H:Y columns one position to the right (H:Z).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before (current) column H.
$ws.Range("H1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("H1").Value = "gross_e5_tu"

# Fill the new column with a straight copy of the gross_e5 (column G) values.
$ws.Range("H2").Value2 = $ws.Range("G2").Value2
$ws.Range("H3").Value2 = $ws.Range("G3").Value2
$ws.Range("H4").Value2 = $ws.Range("G4").Value2
$ws.Range("H5").Value2 = $ws.Range("G5").Value2
$ws.Range("H6").Value2 = $ws.Range("G6").Value2

# The column insert above leaves the Q3:Q6 block (previously the shared
# formula "=O3" living in P3:P6) as individual, non-shared formulas even
# though every cell still holds the same formula text "=P3". Re-assign the
# block as one Formula write so the engine re-collapses it back into a
# single shared-formula group, matching the original authoring shape.
$ws.Range("Q3:Q6").Formula = "=P3"

# Restore the view settings that changed with the edit.
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("G6").Select()
